$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.1423793782173277
$ws.Cells.Item(2, 4).Value = 0.02389231384312041
$ws.Cells.Item(2, 5).Value = 0.1502547449891267
$ws.Cells.Item(2, 6).Value = 0.6662018445176656
$ws.Cells.Item(2, 7).Value = 0.002416400815748054
$ws.Cells.Item(2, 11).Value = 0.5758289050956193
$ws.Cells.Item(2, 13).Value = 0.2662015917416625
$ws.Cells.Item(2, 15).Value = 2.254197210472284

$ws.Cells.Item(3, 2).Value = 0.1329601202308339
$ws.Cells.Item(3, 4).Value = 0.02227087676546802
$ws.Cells.Item(3, 5).Value = 0.1429652403097634
$ws.Cells.Item(3, 6).Value = 0.6630842609392005
$ws.Cells.Item(3, 7).Value = 0.002419150599846586
$ws.Cells.Item(3, 11).Value = 0.50333163287209
$ws.Cells.Item(3, 13).Value = 0.2361896613441345
$ws.Cells.Item(3, 15).Value = 2.258141254488436

$ws.Cells.Item(4, 2).Value = 0.1272476187593412
$ws.Cells.Item(4, 4).Value = 0.02126820484708247
$ws.Cells.Item(4, 5).Value = 0.1386229064383073
$ws.Cells.Item(4, 6).Value = 0.6616748310768727
$ws.Cells.Item(4, 7).Value = 0.002420928902116901
$ws.Cells.Item(4, 11).Value = 0.4586284242983254
$ws.Cells.Item(4, 13).Value = 0.2178025018083716
$ws.Cells.Item(4, 15).Value = 2.26224590360809

$ws.Cells.Item(5, 2).Value = 0.1249377078308243
$ws.Cells.Item(5, 4).Value = 0.02085784934577362
$ws.Cells.Item(5, 5).Value = 0.1368866495342189
$ws.Cells.Item(5, 6).Value = 0.6612272059808149
$ws.Cells.Item(5, 7).Value = 0.002421676255543083
$ws.Cells.Item(5, 11).Value = 0.4403647610978112
$ws.Cells.Item(5, 13).Value = 0.2103197066472049
$ws.Cells.Item(5, 15).Value = 2.264341020683503

$ws.Cells.Item(6, 2).Value = 0.1245552388290747
$ws.Cells.Item(6, 4).Value = 0.0207896047312417
$ws.Cells.Item(6, 5).Value = 0.1366003472668922
$ws.Cells.Item(6, 6).Value = 0.6611605263835187
$ws.Cells.Item(6, 7).Value = 0.002421801724988226
$ws.Cells.Item(6, 11).Value = 0.4373292982072883
$ws.Cells.Item(6, 13).Value = 0.2090778053685156
$ws.Cells.Item(6, 15).Value = 2.264714405234656

$ws.Cells.Item(7, 2).Value = 0.1272163935121711
$ws.Cells.Item(7, 4).Value = 0.02126267773291346
$ws.Cells.Item(7, 5).Value = 0.1385993562905696
$ws.Cells.Item(7, 6).Value = 0.6616682813970627
$ws.Cells.Item(7, 7).Value = 0.00242093888931616
$ws.Cells.Item(7, 11).Value = 0.4583823021661999
$ws.Cells.Item(7, 13).Value = 0.2177015452420292
$ws.Cells.Item(7, 15).Value = 2.262272449668501

$ws.Cells.Item(8, 2).Value = 0.1391169684576568
$ws.Cells.Item(8, 4).Value = 0.02333473255384177
$ws.Cells.Item(8, 5).Value = 0.1477134560887308
$ws.Cells.Item(8, 6).Value = 0.6650220105165943
$ws.Cells.Item(8, 7).Value = 0.002417330319372875
$ws.Cells.Item(8, 11).Value = 0.5508717295973611
$ws.Cells.Item(8, 13).Value = 0.2558451012835974
$ws.Cells.Item(8, 15).Value = 2.255207347972913

$ws.Cells.Item(9, 2).Value = 0.1630123411814282
$ws.Cells.Item(9, 4).Value = 0.02734064688796423
$ws.Cells.Item(9, 5).Value = 0.16665863846017
$ws.Cells.Item(9, 6).Value = 0.675614377461315
$ws.Cells.Item(9, 7).Value = 0.002410964273914238
$ws.Cells.Item(9, 11).Value = 0.7307068231437199
$ws.Cells.Item(9, 13).Value = 0.3309685262020494
$ws.Cells.Item(9, 15).Value = 2.254743804378933

$ws.Cells.Item(10, 2).Value = 0.1809043505654273
$ws.Cells.Item(10, 4).Value = 0.0302477225467328
$ws.Cells.Item(10, 5).Value = 0.1812517653451451
$ws.Cells.Item(10, 6).Value = 0.6858612823085366
$ws.Cells.Item(10, 7).Value = 0.002406715765597791
$ws.Cells.Item(10, 11).Value = 0.861864951420273
$ws.Cells.Item(10, 13).Value = 0.3863710062659322
$ws.Cells.Item(10, 15).Value = 2.262622730070404

$ws.Cells.Item(11, 2).Value = 0.1891160667378955
$ws.Cells.Item(11, 4).Value = 0.03156218389636223
$ws.Cells.Item(11, 5).Value = 0.1880413917265429
$ws.Cells.Item(11, 6).Value = 0.6910617336283309
$ws.Cells.Item(11, 7).Value = 0.002404875139551408
$ws.Cells.Item(11, 11).Value = 0.9213168216000724
$ws.Cells.Item(11, 13).Value = 0.41162337763312
$ws.Cells.Item(11, 15).Value = 2.268004122211238

$ws.Cells.Item(12, 2).Value = 0.1922359518293462
$ws.Cells.Item(12, 4).Value = 0.03205876496600268
$ws.Cells.Item(12, 5).Value = 0.1906345085663332
$ws.Cells.Item(12, 6).Value = 0.693108783477399
$ws.Cells.Item(12, 7).Value = 0.002404191306287341
$ws.Cells.Item(12, 11).Value = 0.9437983902450355
$ws.Cells.Item(12, 13).Value = 0.4211930230845411
$ws.Cells.Item(12, 15).Value = 2.270301268000111

$ws.Cells.Item(13, 2).Value = 0.1915635735948626
$ws.Cells.Item(13, 4).Value = 0.03195187015192857
$ws.Cells.Item(13, 5).Value = 0.190075049781619
$ws.Cells.Item(13, 6).Value = 0.6926644529478381
$ws.Cells.Item(13, 7).Value = 0.002404337997277871
$ws.Cells.Item(13, 11).Value = 0.938958001676184
$ws.Cells.Item(13, 13).Value = 0.4191317119806968
$ws.Cells.Item(13, 15).Value = 2.269794988662198

$ws.Cells.Item(14, 2).Value = 0.1893725362458127
$ws.Cells.Item(14, 4).Value = 0.03160306167278293
$ws.Cells.Item(14, 5).Value = 0.1882542857270835
$ws.Cells.Item(14, 6).Value = 0.6912285860636018
$ws.Cells.Item(14, 7).Value = 0.002404818616292727
$ws.Cells.Item(14, 11).Value = 0.9231670320128842
$ws.Cells.Item(14, 13).Value = 0.4124105354020458
$ws.Cells.Item(14, 15).Value = 2.268187907342082

$ws.Cells.Item(15, 2).Value = 0.1880317978519912
$ws.Cells.Item(15, 4).Value = 0.03138925227855793
$ws.Cells.Item(15, 5).Value = 0.1871418933594526
$ws.Cells.Item(15, 6).Value = 0.6903592085892285
$ws.Cells.Item(15, 7).Value = 0.002405114723122441
$ws.Cells.Item(15, 11).Value = 0.9134904673277617
$ws.Cells.Item(15, 13).Value = 0.4082945480572988
$ws.Cells.Item(15, 15).Value = 2.267237322340691

$ws.Cells.Item(16, 2).Value = 0.1803691430054641
$ws.Cells.Item(16, 4).Value = 0.03016165629312439
$ws.Cells.Item(16, 5).Value = 0.180811114155297
$ws.Cells.Item(16, 6).Value = 0.6855322882013866
$ws.Cells.Item(16, 7).Value = 0.002406837901458167
$ws.Cells.Item(16, 11).Value = 0.8579752685860456
$ws.Cells.Item(16, 13).Value = 0.384721700677062
$ws.Cells.Item(16, 15).Value = 2.262307284422093

$ws.Cells.Item(17, 2).Value = 0.1756868304345858
$ws.Cells.Item(17, 4).Value = 0.02940650074300066
$ws.Cells.Item(17, 5).Value = 0.1769663092248877
$ws.Cells.Item(17, 6).Value = 0.6827093666854864
$ws.Cells.Item(17, 7).Value = 0.002407918544553439
$ws.Cells.Item(17, 11).Value = 0.8238633277364613
$ws.Cells.Item(17, 13).Value = 0.3702732023530402
$ws.Cells.Item(17, 15).Value = 2.259743813440593

$ws.Cells.Item(18, 2).Value = 0.1730005257178107
$ws.Cells.Item(18, 4).Value = 0.02897140566764733
$ws.Cells.Item(18, 5).Value = 0.1747690872077499
$ws.Cells.Item(18, 6).Value = 0.6811364241178808
$ws.Cells.Item(18, 7).Value = 0.002408548768912544
$ws.Cells.Item(18, 11).Value = 0.8042230873016933
$ws.Cells.Item(18, 13).Value = 0.3619674677177542
$ws.Cells.Item(18, 15).Value = 2.25843848339494

$ws.Cells.Item(19, 2).Value = 0.1720921679274028
$ws.Cells.Item(19, 4).Value = 0.02882396210285521
$ws.Cells.Item(19, 5).Value = 0.1740275767377994
$ws.Cells.Item(19, 6).Value = 0.6806125581430393
$ws.Cells.Item(19, 7).Value = 0.002408763642795613
$ws.Cells.Item(19, 11).Value = 0.7975698425095743
$ws.Cells.Item(19, 13).Value = 0.3591560889287635
$ws.Cells.Item(19, 15).Value = 2.25802553884202

$ws.Cells.Item(20, 2).Value = 0.1761845639924502
$ws.Cells.Item(20, 4).Value = 0.02948696613027124
$ws.Cells.Item(20, 5).Value = 0.1773741226271568
$ws.Cells.Item(20, 6).Value = 0.6830046197337367
$ws.Cells.Item(20, 7).Value = 0.002407802611791232
$ws.Cells.Item(20, 11).Value = 0.8274966744705523
$ws.Cells.Item(20, 13).Value = 0.371810787426611
$ws.Cells.Item(20, 15).Value = 2.259999190521313

$ws.Cells.Item(21, 2).Value = 0.1900158190241967
$ws.Cells.Item(21, 4).Value = 0.03170554735014264
$ws.Cells.Item(21, 5).Value = 0.1887884879195951
$ws.Cells.Item(21, 6).Value = 0.6916482229667622
$ws.Cells.Item(21, 7).Value = 0.002404677089947412
$ws.Cells.Item(21, 11).Value = 0.9278060870670686
$ws.Cells.Item(21, 13).Value = 0.4143845146161311
$ws.Cells.Item(21, 15).Value = 2.268652900876305

$ws.Cells.Item(22, 2).Value = 0.1991152327622956
$ws.Cells.Item(22, 4).Value = 0.03314864513956195
$ws.Cells.Item(22, 5).Value = 0.1963769980056966
$ws.Cells.Item(22, 6).Value = 0.6977506212373896
$ws.Cells.Item(22, 7).Value = 0.002402711129361809
$ws.Cells.Item(22, 11).Value = 0.9931796535780109
$ws.Cells.Item(22, 13).Value = 0.4422504121982911
$ws.Cells.Item(22, 15).Value = 2.275820582232143

$ws.Cells.Item(23, 2).Value = 0.1942532729301547
$ws.Cells.Item(23, 4).Value = 0.03237907514817806
$ws.Cells.Item(23, 5).Value = 0.1923150049326168
$ws.Cells.Item(23, 6).Value = 0.694452102671292
$ws.Cells.Item(23, 7).Value = 0.00240375339783859
$ws.Cells.Item(23, 11).Value = 0.9583057509250068
$ws.Cells.Item(23, 13).Value = 0.4273740593485229
$ws.Cells.Item(23, 15).Value = 2.271856415977481

$ws.Cells.Item(24, 2).Value = 0.1759595211373863
$ws.Cells.Item(24, 4).Value = 0.0294505906738749
$ws.Cells.Item(24, 5).Value = 0.1771897090728913
$ws.Cells.Item(24, 6).Value = 0.6828709801224591
$ws.Cells.Item(24, 7).Value = 0.002407854997085129
$ws.Cells.Item(24, 11).Value = 0.8258541281480234
$ws.Cells.Item(24, 13).Value = 0.3711156424185091
$ws.Cells.Item(24, 15).Value = 2.259883209908452

$ws.Cells.Item(25, 2).Value = 0.1564886785348989
$ws.Cells.Item(25, 4).Value = 0.02626319745264283
$ws.Cells.Item(25, 5).Value = 0.1614165539286176
$ws.Cells.Item(25, 6).Value = 0.6723172307177521
$ws.Cells.Item(25, 7).Value = 0.002412610869080404
$ws.Cells.Item(25, 11).Value = 0.6822242962934411
$ws.Cells.Item(25, 13).Value = 0.310609752000282
$ws.Cells.Item(25, 15).Value = 2.264341020683503

